# Powermax105 template: change "data column(s)" -> "database column(s)" in the
# three explanatory header-note cells, and update the sheet view so the window
# is scrolled over to the right-hand columns with U30 as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: "...map to the MATERIAL data column because..." -> "...database column because..."
$a1 = $ws.Range("A1")
$a1.Characters(52, 11).Text = "database column"

# --- G1: rich-text cell "...BASE_FEEDRATE data column by the header attribute in the[ XML][ ]transform."
# Rewriting the whole value (the interop layer doesn't support editing Characters()
# in-place without flattening runs), then re-applying the original per-run font
# formatting (Tahoma 8pt; the " " run stays bold red) so the rich-text layout
# still matches the source as closely as possible.
$g1 = $ws.Range("G1")
$g1.Value() = "This column is mapped to the BASE_FEEDRATE database column by the header attribute in the XML transform."

$g1Run2 = $g1.Characters(90, 4)     # " XML"
$g1Run2.Font.Name = "Tahoma"
$g1Run2.Font.Size = 8

$g1Run3 = $g1.Characters(94, 1)     # " " (bold, red)
$g1Run3.Font.Name = "Tahoma"
$g1Run3.Font.Size = 8
$g1Run3.Font.Bold = $true
$g1Run3.Font.Color = 255

$g1Run4 = $g1.Characters(95, 10)    # "transform."
$g1Run4.Font.Name = "Tahoma"
$g1Run4.Font.Size = 8

# --- U1: "Unmapped data columns and unmapped custom columns..." -> "Unmapped database columns..."
$u1 = $ws.Range("U1")
$u1.Characters(10, 12).Text = "database columns"

# --- Sheet view: scroll the window so column K is the left-most visible
# column, and leave the active selection on U30.
$win = $excel.ActiveWindow
$ws.Range("K1").Select()
$win.ScrollColumn = 11
$win.ScrollRow = 1
$ws.Range("U30").Select()
